$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview: Latest HO Xliff Generate Date for c9315bc4 row.
# This text is shared (same underlying string) with de-de!H2 (Correspond
# Handoff Datetime), so both cells must be updated together.
$wsOverview.Range("G2").Value = "2016-08-15 15:02:22"
$wsDeDe.Range("H2").Value = "2016-08-15 15:02:22"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for c9315bc4 row
$wsZhCn.Range("H2").Value = "2016-08-15 15:02:17"
$wsZhCn.Range("K2").Value = "2016-08-15 15:02:36"

# de-de: Correspond Handback DateTime for c9315bc4 row
$wsDeDe.Range("K2").Value = "2016-08-15 15:02:43"
